$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.939.28"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "2.337.42"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'539.40"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").Value = "'134.42"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = "  +4.57%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +4.35%  "
$ws.Range("D13").Value = "'23.70"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "2.751.33"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "57.882.88"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.338.73"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'10.68"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "'332.55"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "'6.68"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.60"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'62.76"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +7.54%  "
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("D30").Value = "'170.59"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "0.0₃0733"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").Value = "'6.09"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "'18.51"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  +12.73%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.23"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'1.25"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "'1.64"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").Value = "'38.93"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'3.63"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "'285.60"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "'0.0940"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +2.88%  "
$ws.Range("D47").Value = "'0.0502"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0217"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").Value = "'0.381"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").Value = "'17.48"
$ws.Range("E51").Value = "  +1.57%  "
